$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3790618328823996
$ws.Range("C2").Value = 0.1501823978460379
$ws.Range("E2").Value = 0.1573397622411221
$ws.Range("F2").Value = 3.437145544801808
$ws.Range("G2").Value = 0.002551889553826868
$ws.Range("I2").Value = 1.736021053187841
$ws.Range("J2").Value = 0.1648030239805962
$ws.Range("K2").Value = 0.5274573011303687
$ws.Range("M2").Value = 0.3039608317308549

$ws.Range("B3").Value = 0.3565901102324744
$ws.Range("C3").Value = 0.1419382753181395
$ws.Range("E3").Value = 0.1571323646746094
$ws.Range("F3").Value = 3.38593704794944
$ws.Range("G3").Value = 0.002555956128537534
$ws.Range("I3").Value = 1.710023972944001
$ws.Range("J3").Value = 0.1634362458668406
$ws.Range("K3").Value = 0.4969189050714817
$ws.Range("M3").Value = 0.2958907610078541

$ws.Range("B4").Value = 0.3430754657063915
$ws.Range("C4").Value = 0.1369663619482822
$ws.Range("E4").Value = 0.1571287998695432
$ws.Range("F4").Value = 3.355885465935387
$ws.Range("G4").Value = 0.002558583496763259
$ws.Range("I4").Value = 1.694497693740033
$ws.Range("J4").Value = 0.1626174474560074
$ws.Range("K4").Value = 0.4785386308323041
$ws.Range("M4").Value = 0.2911810728126767

$ws.Range("B5").Value = 0.3376392832073236
$ws.Range("C5").Value = 0.1349627664654349
$ws.Range("E5").Value = 0.1571585260922532
$ws.Range("F5").Value = 3.343988188309083
$ws.Range("G5").Value = 0.002559687087062948
$ws.Range("I5").Value = 1.688279673508823
$ws.Range("J5").Value = 0.1622888673035128
$ws.Range("K5").Value = 0.4711414760866148
$ws.Range("M5").Value = 0.2893235508995389

$ws.Range("B6").Value = 0.3367409052326309
$ws.Range("C6").Value = 0.1346314259890136
$ws.Range("E6").Value = 0.157165346847858
$ws.Range("F6").Value = 3.342033712008643
$ws.Range("G6").Value = 0.002559872328650596
$ws.Range("I6").Value = 1.687253746286544
$ws.Range("J6").Value = 0.16223461273065
$ws.Range("K6").Value = 0.4699187931323365
$ws.Range("M6").Value = 0.28901883870288

$ws.Range("B7").Value = 0.3430018634334431
$ws.Range("C7").Value = 0.1369392498354927
$ws.Range("E7").Value = 0.1571290744530849
$ws.Range("F7").Value = 3.355723603001309
$ws.Range("G7").Value = 0.002558598246633457
$ws.Range("I7").Value = 1.694413394524574
$ws.Range("J7").Value = 0.162612995572097
$ws.Range("K7").Value = 0.4784384940866744
$ws.Range("M7").Value = 0.2911557717222024

$ws.Range("B8").Value = 0.3712547955170749
$ws.Range("C8").Value = 0.1473210584389335
$ws.Range("E8").Value = 0.1572425926289007
$ws.Range("F8").Value = 3.419199780909651
$ws.Range("G8").Value = 0.002553264693564184
$ws.Range("I8").Value = 1.726966403447832
$ws.Range("J8").Value = 0.1643274954663951
$ws.Range("K8").Value = 0.516850653302896
$ws.Range("M8").Value = 0.3011273492054087

$ws.Range("B9").Value = 0.4289103896159077
$ws.Range("C9").Value = 0.1684012534249177
$ws.Range("E9").Value = 0.1584453409377922
$ws.Range("F9").Value = 3.554751884821826
$ws.Range("G9").Value = 0.002543835952561643
$ws.Range("I9").Value = 1.794292493323951
$ws.Range("J9").Value = 0.1678537746573028
$ws.Range("K9").Value = 0.5951294658284496
$ws.Range("M9").Value = 0.3226296229012604

$ws.Range("B10").Value = 0.4726557600151011
$ws.Range("C10").Value = 0.184340610055159
$ws.Range("E10").Value = 0.159924420351178
$ws.Range("F10").Value = 3.661165281432687
$ws.Range("G10").Value = 0.002537529875838039
$ws.Range("I10").Value = 1.845931540406013
$ws.Range("J10").Value = 0.1705479265589318
$ws.Range("K10").Value = 0.6544668461505125
$ws.Range("M10").Value = 0.3396191426035813

$ws.Range("B11").Value = 0.4928608256427367
$ws.Range("C11").Value = 0.1916926197440318
$ws.Range("E11").Value = 0.1607261977765226
$ws.Range("F11").Value = 3.711073449213274
$ws.Range("G11").Value = 0.002534794497538796
$ws.Range("I11").Value = 1.869906139521419
$ws.Range("J11").Value = 0.1717967452719193
$ws.Range("K11").Value = 0.6818635670788638
$ws.Range("M11").Value = 0.3476080348901149

$ws.Range("B12").Value = 0.5005559888098787
$ws.Range("C12").Value = 0.1944913582398726
$ws.Range("E12").Value = 0.1610483121111095
$ws.Range("F12").Value = 3.730189146134791
$ws.Range("G12").Value = 0.002533777735522028
$ws.Range("I12").Value = 1.879054980783536
$ws.Range("J12").Value = 0.1722730343401082
$ws.Range("K12").Value = 0.6922964252912038
$ws.Range("M12").Value = 0.3506706951044833

$ws.Range("B13").Value = 0.498896742109423
$ws.Range("C13").Value = 0.1938879439639152
$ws.Range("E13").Value = 0.1609781169530677
$ws.Range("F13").Value = 3.726062593977076
$ws.Range("G13").Value = 0.00253399586707632
$ws.Range("I13").Value = 1.87708148292532
$ws.Range("J13").Value = 0.1721703055098871
$ws.Range("K13").Value = 0.6900469251189065
$ws.Range("M13").Value = 0.3500094310250219

$ws.Range("B14").Value = 0.4934930302147791
$ws.Range("C14").Value = 0.1919225783200318
$ws.Range("E14").Value = 0.1607523278140057
$ws.Range("F14").Value = 3.712641764175032
$ws.Range("G14").Value = 0.002534710466346917
$ws.Range("I14").Value = 1.870657410114887
$ws.Range("J14").Value = 0.1718358616784528
$ws.Range("K14").Value = 0.6827207144798138
$ws.Range("M14").Value = 0.3478592513465557

$ws.Range("B15").Value = 0.4901888256562756
$ws.Range("C15").Value = 0.1907206529774896
$ws.Range("E15").Value = 0.1606164332823212
$ws.Range("F15").Value = 3.704449345394124
$ws.Range("G15").Value = 0.002535150658962449
$ws.Range("I15").Value = 1.866731638064422
$ws.Range("J15").Value = 0.1716314478629428
$ws.Range("K15").Value = 0.678240800508604
$ws.Range("M15").Value = 0.346547080461356

$ws.Range("B16").Value = 0.4713414565440814
$ws.Range("C16").Value = 0.1838621865639993
$ws.Range("E16").Value = 0.1598746131970827
$ws.Range("F16").Value = 3.657933923722027
$ws.Range("G16").Value = 0.002537711312772093
$ws.Range("I16").Value = 1.844374530991345
$ws.Range("J16").Value = 0.1704667849205421
$ws.Range("K16").Value = 0.6526845529528771
$ws.Range("M16").Value = 0.3391022881121373

$ws.Range("B17").Value = 0.459857394090335
$ws.Range("C17").Value = 0.1796807451186737
$ws.Range("E17").Value = 0.1594525238267508
$ws.Range("F17").Value = 3.629783040713079
$ws.Range("G17").Value = 0.002539316257509164
$ws.Range("I17").Value = 1.830783507672493
$ws.Range("J17").Value = 0.1697582867697918
$ws.Range("K17").Value = 0.6371102276656018
$ws.Range("M17").Value = 0.3346018196934892

$ws.Range("B18").Value = 0.4532807756468173
$ws.Range("C18").Value = 0.1772851989672972
$ws.Range("E18").Value = 0.1592218877674725
$ws.Range("F18").Value = 3.613732551047576
$ws.Range("G18").Value = 0.002540251930939658
$ws.Range("I18").Value = 1.823011805110141
$ws.Range("J18").Value = 0.1693529615869025
$ws.Range("K18").Value = 0.6281902849381424
$ws.Range("M18").Value = 0.3320377665614629

$ws.Range("B19").Value = 0.4510589745619313
$ws.Range("C19").Value = 0.1764757366875074
$ws.Range("E19").Value = 0.1591458844864633
$ws.Range("F19").Value = 3.608322349079117
$ws.Range("G19").Value = 0.002540570892718687
$ws.Range("I19").Value = 1.820388235010569
$ws.Range("J19").Value = 0.1692160991427016
$ws.Range("K19").Value = 0.6251766644276984
$ws.Range("M19").Value = 0.3311738293237312

$ws.Range("B20").Value = 0.46107692020027
$ws.Range("C20").Value = 0.1801248815264955
$ws.Range("E20").Value = 0.159496199976175
$ws.Range("F20").Value = 3.632765135436017
$ws.Range("G20").Value = 0.002539144109973135
$ws.Range("I20").Value = 1.832225582319992
$ws.Range("J20").Value = 0.1698334812065951
$ws.Range("K20").Value = 0.6387642068850425
$ws.Range("M20").Value = 0.3350783667859289

$ws.Range("B21").Value = 0.4950790381602417
$ws.Range("C21").Value = 0.1924994538962039
$ws.Range("E21").Value = 0.1608181458401639
$ws.Range("F21").Value = 3.716577903232547
$ws.Range("G21").Value = 0.002534500054566445
$ws.Range("I21").Value = 1.872542406386032
$ws.Range("J21").Value = 0.1719340035984516
$ws.Range("K21").Value = 0.6848710155286426
$ws.Range("M21").Value = 0.3484897950041912

$ws.Range("B22").Value = 0.5175575616494825
$ws.Range("C22").Value = 0.2006726602070614
$ws.Range("E22").Value = 0.1617899194145025
$ws.Range("F22").Value = 3.772617045276689
$ws.Range("G22").Value = 0.00253157598383659
$ws.Range("I22").Value = 1.899301137857236
$ws.Range("J22").Value = 0.1733265892234144
$ws.Range("K22").Value = 0.7153445001086709
$ws.Range("M22").Value = 0.3574731768257635

$ws.Range("B23").Value = 0.505536884795049
$ws.Range("C23").Value = 0.1963025732925985
$ws.Range("E23").Value = 0.1612614149624747
$ws.Range("F23").Value = 3.742592089877149
$ws.Range("G23").Value = 0.002533126483736246
$ws.Range("I23").Value = 1.884981838080577
$ws.Range("J23").Value = 0.1725815152457244
$ws.Range("K23").Value = 0.6990490341567579
$ws.Range("M23").Value = 0.3526586028099103

$ws.Range("B24").Value = 0.4605254923320388
$ws.Range("C24").Value = 0.1799240612409392
$ws.Range("E24").Value = 0.159476416518487
$ws.Range("F24").Value = 3.631416513609679
$ws.Range("G24").Value = 0.00253922189763611
$ws.Range("I24").Value = 1.831573489719858
$ws.Range("J24").Value = 0.1697994795766107
$ws.Range("K24").Value = 0.6380163371679259
$ws.Range("M24").Value = 0.3348628471896902

$ws.Range("B25").Value = 0.413070631904418
$ws.Range("C25").Value = 0.162620038394067
$ws.Range("E25").Value = 0.1580152572623881
$ws.Range("F25").Value = 3.516888399337944
$ws.Range("G25").Value = 0.002546277087632029
$ws.Range("I25").Value = 1.775700531245846
$ws.Range("J25").Value = 0.1668819053058783
$ws.Range("K25").Value = 0.5736342189110815
$ws.Range("M25").Value = 0.3166037144124374
